# Scen_FUEL_PRICE_PROJ.xlsx - add MIN_PINK_HYDROGEN / MIN_GREEN_HYDROGEN
# minimum-price blocks to the "FUEL COST" sheet, mirroring the existing
# IMP_HYDROGEN ("COST" / "IMP" / "HYDROGEN") block that sits in rows 70-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone the layout (styles, borders, number formats, spacer row) of
#        the existing HYDROGEN cost block (rows 70-76) into the two new
#        blocks below it. Copy/paste keeps every column's existing cell
#        style (B/C/D/E/G) exactly as Excel would when a user duplicates a
#        block to extend the table.
$sourceBlock = $ws.Range("B70:G76")
$pinkBlock   = $ws.Range("B77:G83")
$greenBlock  = $ws.Range("B84:G90")

$sourceBlock.Copy($pinkBlock)
$sourceBlock.Copy($greenBlock)

# --- 2. Re-point the "Other_Indexes" column (C) from IMP to MIN for both
#        new blocks - these are minimum-bound rows, not import rows.
$ws.Range("C77:C82").Value = "MIN"
$ws.Range("C84:C89").Value = "MIN"

# --- 3. Process-set names (F) identify the new bound processes.
$ws.Range("F77").Value = "MIN_PINK_HYDROGEN"
$ws.Range("F78").Value = "MIN_PINK_HYDROGEN"
$ws.Range("F79").Value = "MIN_PINK_HYDROGEN"
$ws.Range("F80").Value = "MIN_PINK_HYDROGEN"
$ws.Range("F81").Value = "MIN_PINK_HYDROGEN"
$ws.Range("F82").Value = "MIN_PINK_HYDROGEN"

$ws.Range("F84").Value = "MIN_GREEN_HYDROGEN"
$ws.Range("F85").Value = "MIN_GREEN_HYDROGEN"
$ws.Range("F86").Value = "MIN_GREEN_HYDROGEN"
$ws.Range("F87").Value = "MIN_GREEN_HYDROGEN"
$ws.Range("F88").Value = "MIN_GREEN_HYDROGEN"
$ws.Range("F89").Value = "MIN_GREEN_HYDROGEN"

# Give the pasted-in process-name column its own (imported-looking) plain
# Arial font, distinguishing it from the rest of the row, one named style
# per block - matching the two distinct styles Excel creates when a value
# is pasted in from another workbook.
$pinkStyle = $wb.Styles.Add("Normalny 3")
$pinkStyle.Font.Name = "Arial"
$pinkStyle.Font.Size = 10
$ws.Range("F77:F82").Style = "Normalny 3"

$greenStyle = $wb.Styles.Add("Normalny 4")
$greenStyle.Font.Name = "Arial"
$greenStyle.Font.Size = 10
$ws.Range("F84:F89").Style = "Normalny 4"

# --- 4. Values [MPLN/PJ] (E) for each year, per block.
$ws.Range("E77").Value = 120
$ws.Range("E78").Value = 116.7
$ws.Range("E79").Value = 113.4
$ws.Range("E80").Value = 110.1
$ws.Range("E81").Value = 106.8
$ws.Range("E82").Value = 103.5

$ws.Range("E84").Value = 95
$ws.Range("E85").Value = 93.3
$ws.Range("E86").Value = 91.6
$ws.Range("E87").Value = 89.9
$ws.Range("E88").Value = 88.2
$ws.Range("E89").Value = 86.5

# --- 5. Move the viewport/selection the way the author left the sheet
#        after adding the new rows.
$ws.Range("A59").Select()
$excel.ActiveWindow.ScrollRow = 59
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K82").Select()
